$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price (D) and Volume(1h) (E) columns for rows with changed data ---
# NOTE: a few Price values are plain decimals (e.g. "0.2620", "57.20")
# whose literal text Excel would otherwise mangle by auto-converting the
# cell to a Number and dropping the significant trailing zero. Those are
# written with a leading apostrophe to force a text cell, matching the
# original inline-string formatting exactly.
$ws.Range("D2").Value = "26.095.08"
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("D3").Value = "1.647.61"
$ws.Range("E3").Value = "  -0.70%  "
$ws.Range("D5").Value = "218.47"
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("D6").Value = "0.5203"
$ws.Range("E6").Value = "  -0.60%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("D8").Value = "'0.2620"
$ws.Range("E8").Value = "  -0.67%  "
$ws.Range("D9").Value = "0.06299"
$ws.Range("E9").Value = "  -0.18%  "
$ws.Range("D10").Value = "'20.30"
$ws.Range("E10").Value = "  -1.58%  "
$ws.Range("D11").Value = "0.07645"
$ws.Range("E11").Value = "  -1.80%  "
$ws.Range("D12").Value = "4.581"
$ws.Range("E12").Value = "  +1.81%  "
$ws.Range("D13").Value = "1.637.21"
$ws.Range("E13").Value = "  -1.00%  "
$ws.Range("D14").Value = "1.874.54"
$ws.Range("E14").Value = "  -0.71%  "
$ws.Range("D15").Value = "0.5571"
$ws.Range("E15").Value = "  -1.01%  "
$ws.Range("D16").Value = "0.0₅8122"
$ws.Range("E16").Value = "  +0.83%  "
$ws.Range("E17").Value = "  -0.14%  "
$ws.Range("D18").Value = "26.055.17"
$ws.Range("E18").Value = "  -0.65%  "
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("D20").Value = "4.586"
$ws.Range("E20").Value = "  -2.67%  "
$ws.Range("D21").Value = "194.15"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("E22").Value = "  +2.30%  "
$ws.Range("D23").Value = "5.918"
$ws.Range("E23").Value = "  -1.67%  "
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("D25").Value = "145.02"
$ws.Range("E25").Value = "  -0.26%  "
$ws.Range("D26").Value = "'0.1180"
$ws.Range("E26").Value = "  -2.27%  "
$ws.Range("D27").Value = "7.191"
$ws.Range("E27").Value = "  -0.34%  "
$ws.Range("D30").Value = "0.05427"
$ws.Range("E30").Value = "  -3.87%  "
$ws.Range("D31").Value = "1.267"
$ws.Range("E31").Value = "  -0.62%  "
$ws.Range("E32").Value = "  -1.25%  "
$ws.Range("D33").Value = "3.322"
$ws.Range("E33").Value = "  -1.14%  "
$ws.Range("D34").Value = "1.557"
$ws.Range("E34").Value = "  -2.81%  "
$ws.Range("D35").Value = "2.415"
$ws.Range("E35").Value = "  +0.46%  "
$ws.Range("D36").Value = "2.779"
$ws.Range("E36").Value = "  -0.89%  "
$ws.Range("D37").Value = "0.9419"
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("D38").Value = "'0.5590"
$ws.Range("E38").Value = "  -2.70%  "
$ws.Range("E39").Value = "  -2.09%  "
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("D41").Value = "5.738"
$ws.Range("E41").Value = "  -4.30%  "
$ws.Range("D42").Value = "1.028.05"
$ws.Range("E42").Value = "  -2.26%  "
$ws.Range("D43").Value = "0.8226"
$ws.Range("E43").Value = "  -2.83%  "
$ws.Range("D44").Value = "100.56"
$ws.Range("E44").Value = "  -2.25%  "
$ws.Range("D45").Value = "1.785.76"
$ws.Range("E45").Value = "  -0.74%  "
$ws.Range("E46").Value = "  +7.69%  "
$ws.Range("D47").Value = "'57.20"
$ws.Range("E47").Value = "  -1.59%  "
$ws.Range("D48").Value = "0.9994"
$ws.Range("E48").Value = "  -0.50%  "
$ws.Range("D50").Value = "7.869"
$ws.Range("E50").Value = "  -1.75%  "
$ws.Range("E51").Value = "  -4.20%  "

# --- Rows 28 and 29 swap places (Toncoin moves above EthereumClassic), each
# with its own refreshed price/volume data ---
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "1.534"
$ws.Range("E28").Value = "  +2.40%  "

$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "15.83"
$ws.Range("E29").Value = "  -1.08%  "
